$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A63").Value = "2025/12/05 08:00"
$ws.Range("B63").Value = "-"
$ws.Range("C63").Value = "-"
$ws.Range("D63").Value = "-"
$ws.Range("E63").Value = "-"
$ws.Range("F63").Value = "-"
$ws.Range("G63").Value = "-"
